# Update countries & provincias Spain
# Applies the data refresh represented in the commit diff:
#  - Timestamp in A1 updated
#  - India (row 7), Kazajistan (row 54), Honduras (row 55), Haiti (row 82),
#    Mongolia (row 168) and Butan (row 186) get refreshed case numbers
#  - The Gabon/Guinea/Kirguistan trio (rows 86-88) is re-ranked: Kirguistan
#    jumps to the top of the group with new figures, and Gabon/Guinea keep
#    their previous figures but shift down one rank
#  - The Laos/Santa Lucia/Fiyi/Dominica group (rows 203-206) swaps rank
#    order (values tied, so only the country names move)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "last refreshed" timestamp
$ws.Range("A1").Value = "Datos actualizados a 1 de Julio de 2020 a las 06:13"

# India - row 7
$ws.Range("D7").Value = 347979
$ws.Range("E7").Value = 220403

# Kazajistan - row 54
$ws.Range("B54").Value = 22308
$ws.Range("C54").Value = 489
$ws.Range("D54").Value = 13558
$ws.Range("E54").Value = 8562

# Honduras - row 55
$ws.Range("B55").Value = 19558
$ws.Range("C55").Value = 740
$ws.Range("D55").Value = 2060
$ws.Range("E55").Value = 17001
$ws.Range("G55").Value = 12
$ws.Range("H55").Value = 497

# Haiti - row 82
$ws.Range("D82").Value = 931
$ws.Range("E82").Value = 4939

# Rows 86-88: Kirguistan overtakes Gabon and Guinea in ranking.
# Kirguistan gets fresh figures; Gabon and Guinea keep their previous
# figures but move down one row each.
$ws.Range("A86").Value = "Kirguistan"
$ws.Range("B86").Value = 5506
$ws.Range("C86").Value = 210
$ws.Range("D86").Value = 2443
$ws.Range("E86").Value = 3002
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 4
$ws.Range("H86").Value = 61

$ws.Range("A87").Value = "Gabon"
$ws.Range("B87").Value = 5394
$ws.Range("C87").Value = 0
$ws.Range("D87").Value = 2420
$ws.Range("E87").Value = 2932
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 42

$ws.Range("A88").Value = "Guinea"
$ws.Range("B88").Value = 5391
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 4326
$ws.Range("E88").Value = 1032
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 33

# Mongolia - row 168
$ws.Range("D168").Value = 176
$ws.Range("E168").Value = 44

# Butan - row 186
$ws.Range("D186").Value = 48
$ws.Range("E186").Value = 29

# Rows 203-206: Laos/Santa Lucia and Fiyi/Dominica swap rank (tied values,
# only the country name order changes).
$ws.Range("A203").Value = "Laos"
$ws.Range("A204").Value = "Santa Lucia"
$ws.Range("A205").Value = "Fiyi"
$ws.Range("A206").Value = "Dominica"
